$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 335, shifting existing rows 335-362 down to 336-363.
$ws.Rows.Item(335).Insert()

# Populate the newly inserted row 335 with the new record.
$ws.Cells.Item(335, 1).Value = 5
$ws.Cells.Item(335, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(335, 3).Value = "Maule"
$ws.Cells.Item(335, 4).Value = 44783
$ws.Cells.Item(335, 5).Value = 7
$ws.Cells.Item(335, 6).Value = 100114014
$ws.Cells.Item(335, 7).Value = "Betarraga"
$ws.Cells.Item(335, 8).Value = "Sin especificar"
$ws.Cells.Item(335, 9).Value = "Primera"
$ws.Cells.Item(335, 10).Value = 4000
$ws.Cells.Item(335, 11).Value = 750
$ws.Cells.Item(335, 12).Value = 750
$ws.Cells.Item(335, 13).Value = 750
$ws.Cells.Item(335, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(335, 15).Value = "Región del Maule"
$ws.Cells.Item(335, 16).Value = 150
$ws.Cells.Item(335, 17).Value = 5
$ws.Cells.Item(335, 18).Value = "Hortaliza"
